# Auto-generated edit script: update cryptos list data (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.180.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.34%  '

# Row 3
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.836.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.18%  '

# Row 4
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.09'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.84%  '

# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6623'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.34%  '

# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.06'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.21%  '

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07424'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.28%  '

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2944'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.25%  '

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.09'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.93%  '

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07754'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.21%  '

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.846.01'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.90%  '

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.015'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.14%  '

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6713'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.20%  '

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.20'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.24%  '

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.127'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.17%  '

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008649'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.45%  '

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.175.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.55%  '

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.086.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.27%  '

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '227.52'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.14%  '

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.52'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.14%  '

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.149'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.17%  '

# Row 25
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'BinanceUSD'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.9996'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.26'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.33%  '

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1413'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.93%  '

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.612'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.66%  '

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.02'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.17%  '

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.512'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.88%  '

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.130'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.69%  '

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.054'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.17%  '

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.190'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.07%  '

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05309'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.69%  '

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.872'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.35%  '

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7397'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.93%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.149'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.82%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.647'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.24%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.303.23'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.12%  '

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01793'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.40%  '

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.742'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.14%  '

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.385'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.33%  '

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9191'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.59%  '

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.19%  '

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'XinFinNetwork'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.08265'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +9.28%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.89'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.27%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.988.63'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.49%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5135'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.57%  '

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.16'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.52%  '

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000120'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.83%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.753'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.09%  '

Write-Host "Update complete"